$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" (Strike#) values. Regenerate these values
# (Strike# -> K) for the data rows.
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
